$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4205.1665
$ws.Range("I76").Value = 4966
$ws.Range("J76").Value = 3444.3333
$ws.Range("K76").Value = 4966
$ws.Range("L76").Value = 3444.3333
$ws.Range("M76").Value = -4651
$ws.Range("N76").Value = -4074.3333
$ws.Range("H79").Value = 4205.1665
$ws.Range("I79").Value = 4966
$ws.Range("J79").Value = 3444.3333
$ws.Range("K79").Value = 4966
$ws.Range("L79").Value = 3444.3333
$ws.Range("M79").Value = -3874
$ws.Range("N79").Value = -5628.3333
$ws.Range("H129").Value = 3965
$ws.Range("I129").Value = 1310.9
$ws.Range("K129").Value = 3932.7
$ws.Range("M129").Value = 1067.3
$ws.Range("H137").Value = 8398.362999999999
$ws.Range("I137").Value = 4716.4287
$ws.Range("J137").Value = 11760.131
$ws.Range("K137").Value = 14149.2861
$ws.Range("L137").Value = 35280.393
$ws.Range("M137").Value = -11599.2861
$ws.Range("N137").Value = -40380.393

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 49999
$ws.Range("I34").Value = 49999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 49999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -49728
$ws.Range("N34").ClearContents()
$ws.Range("H61").Value = 7648.087
$ws.Range("I61").Value = 6018.5
$ws.Range("J61").Value = 10703.5625
$ws.Range("K61").Value = 6018.5
$ws.Range("L61").Value = 10703.5625
$ws.Range("M61").Value = -5806.5
$ws.Range("N61").Value = -11127.5625
$ws.Range("H74").Value = 12764.352
$ws.Range("I74").Value = 12636.6
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 12636.6
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -11762.6
$ws.Range("N74").Value = -16748
$ws.Range("H77").Value = 12764.352
$ws.Range("I77").Value = 12636.6
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 63183
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -58815
$ws.Range("N77").Value = -83736
$ws.Range("H122").Value = 2021.174
$ws.Range("I122").Value = 1956.3158
$ws.Range("K122").Value = 5868.9474
$ws.Range("M122").Value = -3418.9474
$ws.Range("H132").Value = 2190.689
$ws.Range("I132").Value = 1704.6923
$ws.Range("K132").Value = 5114.0769
$ws.Range("M132").Value = -2584.0769
$ws.Range("H136").Value = 7648.087
$ws.Range("I136").Value = 6018.5
$ws.Range("J136").Value = 10703.5625
$ws.Range("K136").Value = 18055.5
$ws.Range("L136").Value = 32110.6875
$ws.Range("M136").Value = -15505.5
$ws.Range("N136").Value = -37210.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 25929
$ws.Range("J32").Value = 25929
$ws.Range("L32").Value = 25929
$ws.Range("N32").Value = -26697
$ws.Range("H36").Value = 10040
$ws.Range("I36").Value = 10040
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 10040
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9506
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H134").Value = 8199.885
$ws.Range("I134").Value = 3708.5293
$ws.Range("K134").Value = 11125.5879
$ws.Range("M134").Value = -8590.5879

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 7499.5
$ws.Range("I17").Value = 7499.5
$ws.Range("K17").Value = 7499.5
$ws.Range("M17").Value = -7325.5
$ws.Range("H58").Value = 4909.7407
$ws.Range("I58").Value = 3620.6365
$ws.Range("K58").Value = 3620.6365
$ws.Range("M58").Value = -3417.6365
$ws.Range("H134").Value = 5746.7095
$ws.Range("I134").Value = 3698.2273
$ws.Range("K134").Value = 11094.6819
$ws.Range("M134").Value = -8559.6819
$ws.Range("H136").Value = 4909.7407
$ws.Range("I136").Value = 3620.6365
$ws.Range("K136").Value = 10861.9095
$ws.Range("M136").Value = -8311.9095
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 17499.5
$ws.Range("I5").Value = 13749.75
$ws.Range("J5").Value = 24999
$ws.Range("K5").Value = 13749.75
$ws.Range("L5").Value = 24999
$ws.Range("M5").Value = -13637.75
$ws.Range("N5").Value = -25223
$ws.Range("H11").Value = 549998
$ws.Range("I11").Value = 549998
$ws.Range("K11").Value = 549998
$ws.Range("M11").Value = -549859
$ws.Range("H23").Value = 1999
$ws.Range("J23").Value = 1999
$ws.Range("L23").Value = 1999
$ws.Range("N23").Value = -2445
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H33").Value = 37249.25
$ws.Range("I33").Value = 29998
$ws.Range("K33").Value = 29998
$ws.Range("M33").Value = -29746
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 20000
$ws.Range("K55").Value = 20000
$ws.Range("M55").Value = -19673
$ws.Range("H102").Value = 1216.85
$ws.Range("I102").Value = 1258.75
$ws.Range("K102").Value = 1258.75
$ws.Range("M102").Value = 363.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4126
$ws.Range("I7").Value = 3221.2222
$ws.Range("J7").Value = 5754.6
$ws.Range("K7").Value = 3221.2222
$ws.Range("L7").Value = 5754.6
$ws.Range("M7").Value = -3109.2222
$ws.Range("N7").Value = -5978.6
$ws.Range("H19").Value = 9466.333000000001
$ws.Range("J19").Value = 9466.333000000001
$ws.Range("L19").Value = 9466.333000000001
$ws.Range("N19").Value = -9806.333000000001
$ws.Range("H46").Value = 1782.1072
$ws.Range("I46").Value = 1069.5
$ws.Range("J46").Value = 2178
$ws.Range("K46").Value = 1069.5
$ws.Range("L46").Value = 2178
$ws.Range("M46").Value = -881.5
$ws.Range("N46").Value = -2554
$ws.Range("H126").Value = 4126
$ws.Range("I126").Value = 3221.2222
$ws.Range("J126").Value = 5754.6
$ws.Range("K126").Value = 9663.6666
$ws.Range("L126").Value = 17263.8
$ws.Range("M126").Value = -7193.6666
$ws.Range("N126").Value = -22203.8
$ws.Range("H136").Value = 4005.6956
$ws.Range("I136").Value = 3613.7646
$ws.Range("K136").Value = 10841.2938
$ws.Range("M136").Value = -8291.293799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 2000
$ws.Range("J24").Value = 2000
$ws.Range("L24").Value = 2000
$ws.Range("N24").Value = -2460
$ws.Range("H26").Value = 18255.5
$ws.Range("J26").Value = 20010
$ws.Range("L26").Value = 20010
$ws.Range("N26").Value = -20596
$ws.Range("H29").Value = 14500
$ws.Range("I29").Value = 10000
$ws.Range("J29").Value = 19000
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 19000
$ws.Range("M29").Value = -9710
$ws.Range("N29").Value = -19580
$ws.Range("H34").Value = 62969
$ws.Range("I34").Value = 20026
$ws.Range("J34").Value = 77283.336
$ws.Range("K34").Value = 20026
$ws.Range("L34").Value = 77283.336
$ws.Range("M34").Value = -19823
$ws.Range("N34").Value = -77689.336
$ws.Range("H37").Value = 74997
$ws.Range("J37").Value = 99995
$ws.Range("L37").Value = 99995
$ws.Range("N37").Value = -100401
$ws.Range("H43").Value = 66331.336
$ws.Range("I43").Value = 49499.5
$ws.Range("J43").Value = 99995
$ws.Range("K43").Value = 49499.5
$ws.Range("L43").Value = 99995
$ws.Range("M43").Value = -49350.5
$ws.Range("N43").Value = -100293
$ws.Range("H126").Value = 9861.375
$ws.Range("I126").Value = 6333.8696
$ws.Range("K126").Value = 19001.6088
$ws.Range("M126").Value = -16531.6088
$ws.Range("H136").Value = 3094.3057
$ws.Range("I136").Value = 2623.6924
$ws.Range("J136").Value = 3360.3044
$ws.Range("K136").Value = 7871.0772
$ws.Range("L136").Value = 10080.9132
$ws.Range("M136").Value = -5321.0772
$ws.Range("N136").Value = -15180.9132

